$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-34 down to 31-35
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the new weekly record
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value = "Los Lagos"
$ws.Cells.Item(30, 4).Value = 44522
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = 300000000
$ws.Cells.Item(30, 7).Value = "Espárragos"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 180
$ws.Cells.Item(30, 11).Value = 1800
$ws.Cells.Item(30, 12).Value = 1800
$ws.Cells.Item(30, 13).Value = 1800
$ws.Cells.Item(30, 14).Value = "`$/kilo"
$ws.Cells.Item(30, 15).Value = "Provincia de Linares"
$ws.Cells.Item(30, 16).Value = 1800
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = "Hortaliza"
